$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-13 19:11:47"
$wsZhCn.Range("H2").Value = "2016-03-13 19:12:09"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-13 19:11:51"
$wsDeDe.Range("H2").Value = "2016-03-13 19:12:15"
